# Update the "Förändrad" date column (C) for data rows 2-18:
# change the stored date serial from 45212 to 45221 (2023-10-13 -> 2023-10-22),
# matching the automatic-update commit that refreshed the "last changed" date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # column C ("Förändrad")
    if ($cell.Value2 -eq 45212) {
        $cell.Value2 = 45221
    }
}
